$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before E (existing E:H shift right to F:I).
# The new column E inherits formatting from column D for every row.
$ws.Columns("E").EntireColumn.Insert()

# New column E should be exactly as wide as D/F/G (15.5703125 chars).
$ws.Columns("E").ColumnWidth = 14.666666666666666

# Row 1 header: new E1 = "element2" (brand-new unique string); style already
# copied from D1 by the column insert, and Value2 does not disturb it since
# the style carries no quotePrefix flag.
$ws.Range("E1").Value2 = "element2"

# Rows 2-9 (except row 8, which stays blank): column E duplicates column D's
# value. Because the source style on these rows uses quotePrefix, assigning
# Value2 resets the cell style, so we restore it afterwards by copying D's
# formatting onto E.
$dataRows = 2,3,4,5,6,7,9
foreach ($r in $dataRows) {
    $src = $ws.Range("D$r")
    $dst = $ws.Range("E$r")
    $dst.Value2 = $src.Value2
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
}
$ws.Application.CutCopyMode = $false

# C8 was blank, now holds "Rennes" (matches the other 35000/Rennes rows).
$ws.Range("C8").Value2 = "Rennes"

# Update the active selection shown when the workbook is reopened.
$ws.Range("C14").Select()
